# Add bar codes for Irminger D3
$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")

# Moorings sheet: Mooring OOIBARCODE for the single mooring row
$wsMoorings.Range("A2").Style = "Normal"
$wsMoorings.Range("A2").Value = "OL000606"

# Asset_Cal_Info sheet: Mooring OOIBARCODE (col B) and Sensor OOIBARCODE (col E)
$wsAsset.Range("B2").Value = "OL000606"
$wsAsset.Range("E2").Value = "OL000607"

$wsAsset.Range("B3").Value = "OL000606"
$wsAsset.Range("E3").Value = "OL000607"

$wsAsset.Range("B4").Value = "OL000606"
$wsAsset.Range("E4").Value = "OL000607"

$wsAsset.Range("B5").Value = "OL000606"
$wsAsset.Range("E5").Value = "OL000607"

$wsAsset.Range("B9").Value = "OL000606"
$wsAsset.Range("E9").Value = "OL000608"

$wsAsset.Range("B7").Value = "OL000606"
$wsAsset.Range("E7").Value = "OL000609"

$wsAsset.Range("B11").Value = "OL000606"
$wsAsset.Range("E11").Value = "OL000610"

# Update selections to match final state; selecting a range on a sheet
# also makes that sheet the active one, so doing Moorings first and
# Asset_Cal_Info last leaves Asset_Cal_Info as the active/selected tab.
$wsMoorings.Range("A2").Select()
$wsAsset.Range("H19").Select()
